$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New daily rows appended after the existing data (rows 204-208)
$ws.Cells.Item(204, 1).Value = "25-10-2021"
$ws.Cells.Item(204, 3).Value = 4.77
$ws.Cells.Item(204, 4).Value = 3.46
$ws.Cells.Item(204, 5).Value = 3.44

$ws.Cells.Item(205, 1).Value = "26-10-2021"
$ws.Cells.Item(205, 4).Value = 3.65
$ws.Cells.Item(205, 5).Value = 3.62

$ws.Cells.Item(206, 1).Value = "27-10-2021"
$ws.Cells.Item(206, 4).Value = 3.61
$ws.Cells.Item(206, 5).Value = 3.67

$ws.Cells.Item(207, 1).Value = "28-10-2021"
$ws.Cells.Item(207, 3).Value = 5.28
$ws.Cells.Item(207, 4).Value = 3.71
$ws.Cells.Item(207, 5).Value = 3.55

$ws.Cells.Item(208, 1).Value = "29-10-2021"
$ws.Cells.Item(208, 4).Value = 3.66
$ws.Cells.Item(208, 5).Value = 3.53
